$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows whose Target cluster (column D) is "ECs":
# old rows 8, 9, 10 (MuSCs/ECs, and the ECs-target rows already at 2 and 5)
# We'll delete bottom-up so row indices of earlier rows stay valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Now rows 2-7 hold the remaining 6 combinations (in original relative order):
# 2: ECs/FAPs, 3: ECs/MuSCs, 4: FAPs/FAPs, 5: FAPs/MuSCs, 6: MuSCs/FAPs, 7: MuSCs/MuSCs
# Overwrite every data cell with the new recomputed TPM values.

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1817723333333333
$ws.Range("N2").Value = 0.5453170000000001
$ws.Range("O2").Value = 0.008966262009224884
$ws.Range("P2").Value = 0.008966262009224884
$ws.Range("Q2").Value = 0.1327805693271111
$ws.Range("R2").Value = 1.195025123944
$ws.Range("S2").Value = 0.0002836270760265657
$ws.Range("T2").Value = 0.0002836270760265657

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.09115
$ws.Range("N3").Value = 60.27345
$ws.Range("O3").Value = 0.9910337379907751
$ws.Range("P3").Value = 0.9910337379907752
$ws.Range("Q3").Value = 14.6761296756
$ws.Range("R3").Value = 132.0851670804
$ws.Range("S3").Value = 0.03134907289802703
$ws.Range("T3").Value = 0.03134907289802703

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.409654999999999
$ws.Range("H4").Value = 28.228965
$ws.Range("I4").Value = 0.4074771110502447
$ws.Range("J4").Value = 0.4074771110502448
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1817723333333333
$ws.Range("N4").Value = 0.5453170000000001
$ws.Range("O4").Value = 0.008966262009224884
$ws.Range("P4").Value = 0.008966262009224884
$ws.Range("Q4").Value = 1.710414945211667
$ws.Range("R4").Value = 15.393734506905
$ws.Range("S4").Value = 0.003653546540438518
$ws.Range("T4").Value = 0.003653546540438519

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.409654999999999
$ws.Range("H5").Value = 28.228965
$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.09115
$ws.Range("N5").Value = 60.27345
$ws.Range("O5").Value = 0.9910337379907751
$ws.Range("P5").Value = 0.9910337379907752
$ws.Range("Q5").Value = 189.05079005325
$ws.Range("R5").Value = 1701.45711047925
$ws.Range("S5").Value = 0.4038235645098062
$ws.Range("T5").Value = 0.4038235645098063

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.95234266666667
$ws.Range("H6").Value = 38.857028
$ws.Range("I6").Value = 0.5608901889757016
$ws.Range("J6").Value = 0.5608901889757018
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1817723333333333
$ws.Range("N6").Value = 0.5453170000000001
$ws.Range("O6").Value = 0.008966262009224884
$ws.Range("P6").Value = 0.008966262009224884
$ws.Range("Q6").Value = 2.354377548652889
$ws.Range("R6").Value = 21.189397937876
$ws.Range("S6").Value = 0.0050290883927598
$ws.Range("T6").Value = 0.005029088392759801

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.95234266666667
$ws.Range("H7").Value = 38.857028
$ws.Range("I7").Value = 0.5608901889757016
$ws.Range("J7").Value = 0.5608901889757018
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.09115
$ws.Range("N7").Value = 60.27345
$ws.Range("O7").Value = 0.9910337379907751
$ws.Range("P7").Value = 0.9910337379907752
$ws.Range("Q7").Value = 260.2274593674
$ws.Range("R7").Value = 2342.0471343066
$ws.Range("S7").Value = 0.5558611005829418
$ws.Range("T7").Value = 0.555861100582942
